$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

# Header row (row 1) - all headers are plain text labels
Set-TextCell "A1" "customer_id"
Set-TextCell "B1" "invoice_number"
Set-TextCell "C1" "invoice_date"
Set-TextCell "D1" "status"
Set-TextCell "E1" "is_gst_invoice"
Set-TextCell "F1" "subtotal"
Set-TextCell "G1" "cgst_amount"
Set-TextCell "H1" "sgst_amount"
Set-TextCell "I1" "igst_amount"
Set-TextCell "J1" "total_amount"
Set-TextCell "K1" "notes"
Set-TextCell "L1" "terms"
Set-TextCell "M1" "items"
Set-TextCell "N1" "created_at"
Set-TextCell "O1" "id"

# Row 2
Set-TextCell "A2" "af892bfb-eb9d-40aa-b377-20bb463398bc"
Set-TextCell "B2" "INV-1001"
Set-TextCell "C2" "2025-10-30"
Set-TextCell "D2" "draft"
$ws.Range("E2").Value = $true
$ws.Range("F2").Value = 99
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = 10
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 119
Set-TextCell "K2" "Test invoice"
Set-TextCell "L2" "Net 7"
Set-TextCell "M2" ""
Set-TextCell "N2" "2025-10-30T16:39:16.327Z"
Set-TextCell "O2" "cec08352-98b6-45ac-9791-ef387be342c1"

# Row 3 (note: no M3 cell, matching the source data)
Set-TextCell "A3" "7dd73460-a560-4874-886f-78e863a66d49"
Set-TextCell "B3" "INV-1001"
Set-TextCell "C3" "2025-10-30"
Set-TextCell "D3" "draft"
$ws.Range("E3").Value = $true
$ws.Range("F3").Value = 99
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 10
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 119
Set-TextCell "K3" "Test invoice"
Set-TextCell "L3" "Net 7"
Set-TextCell "N3" "2025-10-30T16:49:53.497Z"
Set-TextCell "O3" "9194b2d3-dacb-4b03-acbc-7e984c1d9afa"
